# Generate Report for Archive
#
# Two files (53b72932-d92f-4aba-bf3e-42757ad1db87.md and
# 60bdeb2e-d7d1-4f55-ab26-b15f75c89bb8.md) have moved on from "Ready for
# handoff" and are now "In Translation". Update their status on every
# sheet that tracks it: the Overview sheet (zh-cn / de-de columns) and
# the per-locale detail sheets (Status column).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 -> 53b72932-..., row 4 -> 60bdeb2e-...
# Columns E (zh-cn) and F (de-de) hold the per-locale status.
foreach ($row in 3, 4) {
    if ($overview.Range("E$row").Value2 -eq $oldStatus) {
        $overview.Range("E$row").Value = $newStatus
    }
    if ($overview.Range("F$row").Value2 -eq $oldStatus) {
        $overview.Range("F$row").Value = $newStatus
    }
}

# Per-locale sheets: rows 3 and 4 hold the same two files, column C is Status.
foreach ($sheet in $zhcn, $dede) {
    foreach ($row in 3, 4) {
        if ($sheet.Range("C$row").Value2 -eq $oldStatus) {
            $sheet.Range("C$row").Value = $newStatus
        }
    }
}
